$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.412.12'
$ws.Range('E2').Value = '  -4.87%  '
$ws.Range('D3').Value = '3.263.61'
$ws.Range('E3').Value = '  -7.68%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'588.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.63%  '
$ws.Range('D6').Value = "'153.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.15%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '3.254.69'
$ws.Range('E8').Value = '  -7.88%  '
$ws.Range('D9').Value = "'0.546"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.74%  '
$ws.Range('D10').Value = "'0.173"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -13.50%  '
$ws.Range('D11').Value = "'6.81"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.81%  '
$ws.Range('D12').Value = "'0.509"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.76%  '
$ws.Range('D13').Value = "'38.67"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -16.98%  '
$ws.Range('D14').Value = "'0.0000245"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.23%  '
$ws.Range('D15').Value = '3.781.06'
$ws.Range('E15').Value = '  -7.77%  '
$ws.Range('D16').Value = '67.424.24'
$ws.Range('E16').Value = '  -4.88%  '
$ws.Range('D17').Value = "'550.51"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -9.89%  '
$ws.Range('D18').Value = '3.265.42'
$ws.Range('E18').Value = '  -7.52%  '
$ws.Range('D19').Value = "'7.31"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -13.27%  '
$ws.Range('E20').Value = '  -5.88%  '
$ws.Range('D21').Value = "'15.26"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.06%  '
$ws.Range('D22').Value = "'0.771"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -12.89%  '
$ws.Range('D23').Value = "'7.92"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -12.89%  '
$ws.Range('D24').Value = "'86.00"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -12.83%  '
$ws.Range('D25').Value = "'13.67"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -12.51%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -14.57%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'29.61"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -12.58%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = "'8.10"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.80%  '
$ws.Range('E30').Value = '  -16.66%  '
$ws.Range('D31').Value = "'2.70"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.93%  '
$ws.Range('D32').Value = "'1.16"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.85%  '
$ws.Range('D33').Value = "'552.99"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -11.96%  '
$ws.Range('D34').Value = "'6.70"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -17.76%  '
$ws.Range('D35').Value = "'5.79"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -15.04%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = "'0.0451"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.95%  '
$ws.Range('D38').Value = "'53.64"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.87%  '
$ws.Range('D39').Value = "'0.0864"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.47%  '
$ws.Range('D40').Value = "'9.27"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.36%  '
$ws.Range('E41').Value = '  -12.33%  '
$ws.Range('D42').Value = '2.955.24'
$ws.Range('E42').Value = '  -12.09%  '
$ws.Range('D43').Value = "'2.66"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -23.06%  '
$ws.Range('D44').Value = "'0.264"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -15.28%  '
$ws.Range('D45').Value = '0.0₃0587'
$ws.Range('E45').Value = '  -20.03%  '
$ws.Range('D46').Value = "'26.65"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -16.84%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = "'2.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -15.22%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = "'2.38"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -19.94%  '
$ws.Range('D50').Value = "'127.44"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.48%  '
$ws.Range('D51').Value = "'0.115"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -11.98%  '
